# Horarios Línea 141 - actualización 04:56:30 (20260126)
# Updates the three schedule sheets (LP1912, LP1912-215, 6203-6173) with a
# refreshed scrape: new "Última actualización" timestamp, updated row counts,
# recalculated "Minutos" countdowns, a couple of re-ordered/re-labelled stops,
# and new trailing rows appended to each sheet.

$wb = $excel.ActiveWorkbook

function Set-Row {
    param(
        $ws,
        [int]$row,
        $values   # ordered hashtable-like array of [col]=value pairs, col in A..E
    )
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 04:56:30"
$ws1.Range("A3").Value = "Total filas: 20"

Set-Row $ws1 9  @{ A = "04:56:30"; D = 20 }
Set-Row $ws1 10 @{ A = "04:56:30"; D = 26 }
Set-Row $ws1 11 @{ A = "04:56:30"; D = 38 }
Set-Row $ws1 13 @{ A = "04:56:30"; D = 50 }
Set-Row $ws1 14 @{ A = "04:56:30"; D = 58 }
Set-Row $ws1 16 @{ A = "04:56:30"; D = 75 }
Set-Row $ws1 17 @{ A = "04:56:30"; D = 78 }
Set-Row $ws1 18 @{ A = "04:56:30"; B = "06:18"; C = "16_SANTA ANA"; D = 82 }
Set-Row $ws1 19 @{ A = "04:56:30"; B = "06:21"; C = "26_HERNANDEZ"; D = 85 }
Set-Row $ws1 20 @{ B = "06:24"; C = "16_SANTA ANA"; D = 99 }
Set-Row $ws1 21 @{ A = "04:56:30"; B = "06:27"; C = "23_HERNANDEZ"; D = 91 }
Set-Row $ws1 22 @{ A = "04:56:30"; B = "06:29"; C = "86_EST CHICA-ESC AGRARIA"; D = 93 }
Set-Row $ws1 23 @{ A = "04:56:30"; B = "06:31"; C = "16_SANTA ANA"; D = 95 }
Set-Row $ws1 24 @{ A = "04:56:30"; B = "06:44"; C = "225_C ROCA-H SUR"; D = 108; E = "LP1912" }
Set-Row $ws1 25 @{ A = "04:56:30"; B = "06:46"; C = "215C_EL PATO"; D = 110; E = "LP1912" }

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 04:56:30"
$ws2.Range("A3").Value = "Total filas: 6"

Set-Row $ws2 8  @{ A = "04:56:30"; D = 38 }
Set-Row $ws2 10 @{ A = "04:56:30"; D = 75 }
Set-Row $ws2 11 @{ A = "04:56:30"; B = "06:46"; C = "215C_EL PATO"; D = 110; E = "LP1912" }

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 04:56:30"
$ws3.Range("A3").Value = "Total filas: 6"

Set-Row $ws3 7  @{ A = "04:56:30"; D = 48 }
Set-Row $ws3 9  @{ A = "04:56:30"; D = 73 }
Set-Row $ws3 11 @{ A = "04:56:30"; B = "06:33"; C = "215C_LA PLATA"; D = 97; E = "L6203" }
